$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All target cells hold plain text in the source data (inline strings in the
# original OOXML), including values that look numeric (e.g. "0.997", "16.60").
# A bare Value assignment lets Excel coerce numeric-looking text to a real
# number (dropping formatting like trailing zeros), so every text write below
# uses the classic leading-apostrophe text-prefix, then ClearFormats() to drop
# the quotePrefix cell style Excel applies for that prefix -- leaving the cell
# as plain text with no style index, matching the original formatting.

# Row 2
$ws.Range("D2").Value = "'27.831.45"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +3.02%  "
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").Value = "'1.724.03"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +3.06%  "
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  -0.27%  "
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").Value = "'216.87"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +0.83%  "
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("D6").Value = "'0.522"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +1.11%  "
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  -0.28%  "
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("D8").Value = "'23.88"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +11.44%  "
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("E9").Value = "'  +4.32%  "
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("D10").Value = "'0.0629"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +1.11%  "
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("D11").Value = "'0.0899"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +1.15%  "
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("D12").Value = "'1.968.17"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +3.10%  "
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("D13").Value = "'1.724.66"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +2.25%  "
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("D14").Value = "'4.23"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +2.96%  "
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("D15").Value = "'0.564"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +5.62%  "
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("D16").Value = "'67.95"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +2.63%  "
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("D17").Value = "'27.825.39"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +3.06%  "
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("D18").Value = "'240.79"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +2.37%  "
$ws.Range("E18").ClearFormats()

# Row 19
$ws.Range("D19").Value = "'8.04"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  -2.03%  "
$ws.Range("E19").ClearFormats()

# Row 20
$ws.Range("D20").Value = "'0.0₃0748"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +1.68%  "
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("D21").Value = "'0.997"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -0.39%  "
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("E22").Value = "'  +3.61%  "
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("D23").Value = "'9.68"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +4.48%  "
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("E24").Value = "'  +0.32%  "
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("D25").Value = "'148.43"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +0.86%  "
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("D26").Value = "'7.54"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +3.95%  "
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("D27").Value = "'16.60"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +1.02%  "
$ws.Range("E27").ClearFormats()

# Row 28
$ws.Range("E28").Value = "'  +1.24%  "
$ws.Range("E28").ClearFormats()

# Row 29
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -0.44%  "
$ws.Range("E29").ClearFormats()

# Row 30
$ws.Range("D30").Value = "'0.0505"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +1.34%  "
$ws.Range("E30").ClearFormats()

# Row 31
$ws.Range("D31").Value = "'1.18"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +1.09%  "
$ws.Range("E31").ClearFormats()

# Row 32
$ws.Range("E32").Value = "'  +2.04%  "
$ws.Range("E32").ClearFormats()

# Row 33
$ws.Range("D33").Value = "'3.30"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +4.08%  "
$ws.Range("E33").ClearFormats()

# Row 34
$ws.Range("D34").Value = "'1.467.33"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -4.65%  "
$ws.Range("E34").ClearFormats()

# Row 35
$ws.Range("D35").Value = "'1.68"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -2.03%  "
$ws.Range("E35").ClearFormats()

# Row 36
$ws.Range("D36").Value = "'0.964"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +5.70%  "
$ws.Range("E36").ClearFormats()

# Row 37
$ws.Range("D37").Value = "'0.614"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +3.85%  "
$ws.Range("E37").ClearFormats()

# Row 38
$ws.Range("D38").Value = "'2.40"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +0.41%  "
$ws.Range("E38").ClearFormats()

# Row 39
$ws.Range("E39").Value = "'  -0.01%  "
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("E40").Value = "'  +2.61%  "
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("D41").Value = "'71.65"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +5.94%  "
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("D42").Value = "'5.86"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +5.74%  "
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("D43").Value = "'0.997"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  -0.38%  "
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("B44").Value = "'RocketPoolETH"
$ws.Range("B44").ClearFormats()
$ws.Range("C44").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("C44").ClearFormats()
$ws.Range("D44").Value = "'1.870.96"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +2.98%  "
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("B45").Value = "'MXToken"
$ws.Range("B45").ClearFormats()
$ws.Range("C45").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C45").ClearFormats()
$ws.Range("D45").Value = "'2.28"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +1.83%  "
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("D46").Value = "'0.790"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +1.28%  "
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("E47").Value = "'  +9.45%  "
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("D48").Value = "'91.66"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +1.25%  "
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("D49").Value = "'0.0₆0109"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +3.02%  "
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("B50").Value = "'Algorand"
$ws.Range("B50").ClearFormats()
$ws.Range("C50").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C50").ClearFormats()
$ws.Range("D50").Value = "'0.106"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +2.24%  "
$ws.Range("E50").ClearFormats()

# Row 51
$ws.Range("B51").Value = "'EnergySwap"
$ws.Range("B51").ClearFormats()
$ws.Range("C51").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C51").ClearFormats()
$ws.Range("D51").Value = "'8.28"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +3.29%  "
$ws.Range("E51").ClearFormats()
